$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally listed years 2000-2013 (rows 2-15) with their B-column
# values. The update drops the 2000-2009 rows (old rows 2-11) entirely, so the
# previously-last four years (2010-2013, old rows 12-15) shift up to become
# rows 2-5 and the used range shrinks from A1:B15 to A1:B5.
$ws.Rows("2:11").Delete()
